$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C ("Schaden" and "QS"),
# pushing the old C:G ("Dauer"..."Wirkdauer") to E:I.
$ws.Range("C:D").Insert()

# The inserted columns lose their custom width; match column B's width
# (all of B:D share the same 24-character width).
$ws.Range("C:D").ColumnWidth = $ws.Range("B1").ColumnWidth

# Fill in the new Schaden column
$ws.Range("C1").Value = "Schaden"
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = "2 w6"

# Fill in the new QS column
$ws.Range("D1").Value = "QS"
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = "x2"

# Rename the (now shifted) Dauer column to Zauberdauer and update its
# value labels to the "N Aktion(en)" phrasing
$ws.Range("E1").Value = "Zauberdauer"
$ws.Range("E2").Value = "4 Aktion(en)"
$ws.Range("E3").Value = "1 Aktion(en)"
$ws.Range("E4").Value = "2 Aktion(en)"

# Fix capitalisation of the Kosten value for Flim Flam
$ws.Range("F3").Value = "2 AsP"

# Leave the cursor on the top-left cell (matches the saved selection state)
$ws.Range("A1").Select() | Out-Null
